$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.4
$ws.Range("D3").Value = 5.2
$ws.Range("C4").Value = 10.2
$ws.Range("B5").Value = 9.1
$ws.Range("C5").Value = 9.199999999999999
$ws.Range("B6").Value = 11.7
$ws.Range("C6").Value = 11.8
$ws.Range("D6").Value = 16.1
$ws.Range("D7").Value = 3.7
$ws.Range("C9").Value = 6.3
$ws.Range("B10").Value = 12.7
$ws.Range("C10").Value = 16.6
$ws.Range("D10").Value = 14.5
$ws.Range("B11").Value = 29.6
$ws.Range("D11").Value = 23.8
$ws.Range("C12").Value = 2.9
